$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing "length start end" numeric suffix from the libelle
# (column F) descriptions - dates are now typed as "c" instead of "dmy"
# and the position/length info no longer needs to be repeated in the label.
$ws.Range("F3").Value  = "N° FINESS du fichier d'entrée"
$ws.Range("F4").Value  = "Type de prestation"
$ws.Range("F5").Value  = "Année période"
$ws.Range("F6").Value  = "N° période (mois)"
$ws.Range("F7").Value  = "N° d'index du RSA"
$ws.Range("F8").Value  = "Mois du séjour"
$ws.Range("F9").Value  = "Année du séjour"
$ws.Range("F10").Value = "Nombre d'IVG antérieures"
$ws.Range("F11").Value = "Année de la dernière IVG"
$ws.Range("F12").Value = "Nombre de naissances vivantes antérieures"
$ws.Range("F13").Value = "Filler"

# Move the active selection to F14
$ws.Range("F14").Select()

# Default column width changed slightly (10.59375 -> 10.5859375)
$ws.StandardWidth = 10.5859375
